$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update progress dates (Cap nhat tien do)
$ws.Range("G10").Value = (Get-Date -Year 2018 -Month 10 -Day 27 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("H10").Value = (Get-Date -Year 2018 -Month 10 -Day 28 -Hour 0 -Minute 0 -Second 0).Date

$ws.Range("F11").Value = (Get-Date -Year 2018 -Month 10 -Day 31 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("G11").Value = (Get-Date -Year 2018 -Month 10 -Day 28 -Hour 0 -Minute 0 -Second 0).Date

# Move the active selection to G11, matching the saved cursor position
$ws.Range("G11").Select()
